$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
#   - OVERALL row (row 2): trade count/win-rate/PnL stats updated
#   - leadlag row (row 3): trade count/win-rate/PnL stats updated
#
# Note: percentage-looking text (e.g. "65.9%") and plain decimal text (e.g.
# "2.21") would otherwise be auto-coerced into numbers/percentages by the
# COM layer, same as typing them into Excel would. A leading apostrophe
# forces them to stay literal text, matching the source workbook where
# these columns are stored as inline strings, not numbers.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("C2").Value = 41
$wsSummary.Range("D2").Value = "'65.9%"
$wsSummary.Range("E2").Value = "'+8.9095%"
$wsSummary.Range("F2").Value = "'+0.2173%"

$wsSummary.Range("C3").Value = 50
$wsSummary.Range("D3").Value = "'36.0%"
$wsSummary.Range("E3").Value = "'+4.7079%"
$wsSummary.Range("F3").Value = "'+0.0942%"

# ---------------------------------------------------------------------------
# Sheet: leadlag
#   - Row 31 (Trade #41): the open trade is now closed (time-exit @5min)
#   - New row 52 (Trade #63): a freshly opened DOWN trade
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Cells.Item(31, 7).Value = 68768.717112
$wsLeadlag.Cells.Item(31, 8).Value = "CLOSED"
$wsLeadlag.Cells.Item(31, 9).Value = 0.07870000000000001
$wsLeadlag.Cells.Item(31, 10).Value = 0.79
$wsLeadlag.Cells.Item(31, 13).Value = "time_exit_5min"
$wsLeadlag.Cells.Item(31, 14).Value = 5

$wsLeadlag.Cells.Item(52, 1).Value = 63
$wsLeadlag.Cells.Item(52, 2).Value = "'2026-02-16"
$wsLeadlag.Cells.Item(52, 3).Value = "21:34:21"
$wsLeadlag.Cells.Item(52, 4).Value = "leadlag"
$wsLeadlag.Cells.Item(52, 5).Value = "DOWN"
$wsLeadlag.Cells.Item(52, 6).Value = 68716.735
$wsLeadlag.Cells.Item(52, 8).Value = "OPEN"
$wsLeadlag.Cells.Item(52, 9).Value = 0
$wsLeadlag.Cells.Item(52, 10).Value = 0
$wsLeadlag.Cells.Item(52, 11).Value = 0.7429
$wsLeadlag.Cells.Item(52, 12).Value = "Binance leading with -0.074% move"
$wsLeadlag.Cells.Item(52, 14).Value = 0

# ---------------------------------------------------------------------------
# Sheet: All Trades
#   - New row 42 mirrors leadlag trade #41 now that it has closed
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(42, 1).Value = 41
$wsAll.Cells.Item(42, 2).Value = "'2026-02-16"
$wsAll.Cells.Item(42, 3).Value = "21:29:18"
$wsAll.Cells.Item(42, 4).Value = "leadlag"
$wsAll.Cells.Item(42, 5).Value = "UP"
$wsAll.Cells.Item(42, 6).Value = 68714.66
$wsAll.Cells.Item(42, 7).Value = 68768.717112
$wsAll.Cells.Item(42, 8).Value = "CLOSED"
$wsAll.Cells.Item(42, 9).Value = 0.07870000000000001
$wsAll.Cells.Item(42, 10).Value = 0.79
$wsAll.Cells.Item(42, 11).Value = 0.75
$wsAll.Cells.Item(42, 12).Value = "Binance leading with 0.114% move"
$wsAll.Cells.Item(42, 13).Value = "time_exit_5min"
$wsAll.Cells.Item(42, 14).Value = 5

# ---------------------------------------------------------------------------
# Sheet: Comparison
#   - leadlag row (row 2): trade count/win-rate/profit-factor/ratios updated
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

$wsComparison.Range("B2").Value = 50
$wsComparison.Range("C2").Value = "'36.0%"
$wsComparison.Range("D2").Value = "'2.21"
$wsComparison.Range("E2").Value = "'+0.4773%"
$wsComparison.Range("G2").Value = "'1.48"
